$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 0  # was 5333.3335
$ws.Range("J26").Value = 0  # was 5333.3335
$ws.Range("L26").ClearContents()  # was 5333.3335
$ws.Range("N26").Value = 0  # was -6021.3335
$ws.Range("H64").Value = 3773.2222  # was 3871.0908
$ws.Range("I64").Value = 3770.5  # was 3884.6086
$ws.Range("J64").Value = 3780.3  # was 3840
$ws.Range("K64").Value = 3770.5  # was 3884.6086
$ws.Range("L64").Value = 3780.3  # was 3840
$ws.Range("M64").Value = -3522.5  # was -3636.6086
$ws.Range("N64").Value = -4276.3  # was -4336
$ws.Range("H67").Value = 3773.2222  # was 3871.0908
$ws.Range("I67").Value = 3770.5  # was 3884.6086
$ws.Range("J67").Value = 3780.3  # was 3840
$ws.Range("K67").Value = 3770.5  # was 3884.6086
$ws.Range("L67").Value = 3780.3  # was 3840
$ws.Range("M67").Value = -2912.5  # was -3026.6086
$ws.Range("N67").Value = -5496.3  # was -5556
$ws.Range("H74").Value = 4700  # was 4327.7856
$ws.Range("I74").Value = 4040.3  # was 4170
$ws.Range("J74").Value = 5799.5  # was 4722.25
$ws.Range("K74").Value = 4040.3  # was 4170
$ws.Range("L74").Value = 5799.5  # was 4722.25
$ws.Range("M74").Value = -3104.3  # was -3234
$ws.Range("N74").Value = -7671.5  # was -6594.25
$ws.Range("H77").Value = 4700  # was 4327.7856
$ws.Range("I77").Value = 4040.3  # was 4170
$ws.Range("J77").Value = 5799.5  # was 4722.25
$ws.Range("K77").Value = 20201.5  # was 20850
$ws.Range("L77").Value = 28997.5  # was 23611.25
$ws.Range("M77").Value = -15521.5  # was -16170
$ws.Range("N77").Value = -38357.5  # was -32971.25
$ws.Range("H112").Value = 4680.373  # was 5061.6416
$ws.Range("I112").Value = 495  # was 663.3333
$ws.Range("J112").Value = 4827.228  # was 5325.54
$ws.Range("K112").Value = 1485  # was 1989.9999
$ws.Range("L112").Value = 14481.684  # was 15976.62
$ws.Range("M112").Value = -377  # was -881.9999
$ws.Range("N112").Value = -16697.684  # was -18192.62
$ws.Range("H129").Value = 890.1236  # was 889.4269399999999
$ws.Range("I129").Value = 463.58334  # was 451.08334
$ws.Range("J129").Value = 956.5974  # was 957.74023
$ws.Range("K129").Value = 1390.75002  # was 1353.25002
$ws.Range("L129").Value = 2869.7922  # was 2873.22069
$ws.Range("M129").Value = 3609.24998  # was 3646.74998
$ws.Range("N129").Value = -12869.7922  # was -12873.22069

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5479.2646  # was 5211.176
$ws.Range("I32").Value = 4554.1963  # was 4192.9834
$ws.Range("J32").Value = 9796.25  # was 9988.846
$ws.Range("K32").Value = 4554.1963  # was 4192.9834
$ws.Range("L32").Value = 9796.25  # was 9988.846
$ws.Range("M32").Value = -4267.1963  # was -3905.9834
$ws.Range("N32").Value = -10370.25  # was -10562.846
$ws.Range("H125").Value = 98000  # was 61000
$ws.Range("J125").Value = 98000  # was 61000
$ws.Range("L125").Value = 98000  # was 61000
$ws.Range("N125").Value = -107840  # was -70840

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 19609736  # was 13335065
$ws.Range("I86").Value = 27779476  # was 18520084
$ws.Range("J86").Value = 2361.4  # was 2160.1428
$ws.Range("K86").Value = 27779476  # was 18520084
$ws.Range("L86").Value = 2361.4  # was 2160.1428
$ws.Range("M86").Value = -27778353  # was -18518961
$ws.Range("N86").Value = -4607.4  # was -4406.1428
$ws.Range("H89").Value = 19609736  # was 13335065
$ws.Range("I89").Value = 27779476  # was 18520084
$ws.Range("J89").Value = 2361.4  # was 2160.1428
$ws.Range("K89").Value = 138897380  # was 92600420
$ws.Range("L89").Value = 11807  # was 10800.714
$ws.Range("M89").Value = -138891764  # was -92594804
$ws.Range("N89").Value = -23039  # was -22032.714
$ws.Range("H99").Value = 50001076  # was 45455600
$ws.Range("I99").Value = 66667680  # was 62501004
$ws.Range("J99").Value = 1262.2  # was 1185.1666
$ws.Range("K99").Value = 66667680  # was 62501004
$ws.Range("L99").Value = 1262.2  # was 1185.1666
$ws.Range("M99").Value = -66666182  # was -62499506
$ws.Range("N99").Value = -4258.2  # was -4181.1666
$ws.Range("H105").Value = 21236.727  # was 17665.54
$ws.Range("I105").Value = 51065  # was 29658.285
$ws.Range("J105").Value = 4192  # was 3674
$ws.Range("K105").Value = 51065  # was 29658.285
$ws.Range("L105").Value = 4192  # was 3674
$ws.Range("M105").Value = -49318  # was -27911.285
$ws.Range("N105").Value = -7686  # was -7168

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1876.5  # was 1801.75
$ws.Range("I68").Value = 453  # was 442.8
$ws.Range("J68").Value = 3300  # was 4066.6667
$ws.Range("K68").Value = 1359  # was 1328.4
$ws.Range("L68").Value = 9900  # was 12200.0001
$ws.Range("M68").Value = -548  # was -517.4000000000001
$ws.Range("N68").Value = -11522  # was -13822.0001
$ws.Range("H71").Value = 1876.5  # was 1801.75
$ws.Range("I71").Value = 453  # was 442.8
$ws.Range("J71").Value = 3300  # was 4066.6667
$ws.Range("K71").Value = 4077  # was 3985.2
$ws.Range("L71").Value = 29700  # was 36600.0003
$ws.Range("M71").Value = -21  # was 70.79999999999973
$ws.Range("N71").Value = -37812  # was -44712.0003
$ws.Range("H92").Value = 582.8570999999999  # was 700
$ws.Range("I92").Value = 520  # was 0
$ws.Range("J92").Value = 666.6667  # was 700
$ws.Range("K92").Value = 1560  # was 0
$ws.Range("L92").Value = 2000.0001  # was 2100
$ws.Range("M92").Value = -312  # was None
$ws.Range("N92").Value = -4496.0001  # was -4596
$ws.Range("H116").Value = 2425  # was 0
$ws.Range("I116").Value = 850  # was 0
$ws.Range("J116").Value = 4000  # was 0
$ws.Range("K116").Value = 2550  # was 0
$ws.Range("L116").Value = 12000  # was 0
$ws.Range("M116").Value = 892  # was None
$ws.Range("N116").Value = -18884  # was None

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2108.9744  # was 2094.5945
$ws.Range("I80").Value = 2078.6667  # was 2069.4443
$ws.Range("J80").Value = 2866.6667  # was 3000
$ws.Range("K80").Value = 2078.6667  # was 2069.4443
$ws.Range("L80").Value = 2866.6667  # was 3000
$ws.Range("M80").Value = -1080.6667  # was -1071.4443
$ws.Range("N80").Value = -4862.6667  # was -4996
$ws.Range("H83").Value = 2108.9744  # was 2094.5945
$ws.Range("I83").Value = 2078.6667  # was 2069.4443
$ws.Range("J83").Value = 2866.6667  # was 3000
$ws.Range("K83").Value = 10393.3335  # was 10347.2215
$ws.Range("L83").Value = 14333.3335  # was 15000
$ws.Range("M83").Value = -5401.333500000001  # was -5355.2215
$ws.Range("N83").Value = -24317.3335  # was -24984

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 111113560  # was 40002110
$ws.Range("I68").Value = 1475  # was 1691.091
$ws.Range("J68").Value = 142859870  # was 71431016
$ws.Range("K68").Value = 1475  # was 1691.091
$ws.Range("L68").Value = 142859870  # was 71431016
$ws.Range("M68").Value = -726  # was -942.0909999999999
$ws.Range("N68").Value = -142861368  # was -71432514
$ws.Range("H71").Value = 111113560  # was 40002110
$ws.Range("I71").Value = 1475  # was 1691.091
$ws.Range("J71").Value = 142859870  # was 71431016
$ws.Range("K71").Value = 7375  # was 8455.455
$ws.Range("L71").Value = 714299350  # was 357155080
$ws.Range("M71").Value = -3631  # was -4711.455
$ws.Range("N71").Value = -714306838  # was -357162568
$ws.Range("H82").Value = 94705.63  # was 87021.836
$ws.Range("I82").Value = 1565.5  # was 1512.4
$ws.Range("J82").Value = 147928.58  # was 148100
$ws.Range("K82").Value = 1565.5  # was 1512.4
$ws.Range("L82").Value = 147928.58  # was 148100
$ws.Range("M82").Value = -1204.5  # was -1151.4
$ws.Range("N82").Value = -148650.58  # was -148822
$ws.Range("H85").Value = 94705.63  # was 87021.836
$ws.Range("I85").Value = 1565.5  # was 1512.4
$ws.Range("J85").Value = 147928.58  # was 148100
$ws.Range("K85").Value = 1565.5  # was 1512.4
$ws.Range("L85").Value = 147928.58  # was 148100
$ws.Range("M85").Value = -317.5  # was -264.4000000000001
$ws.Range("N85").Value = -150424.58  # was -150596

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 45000  # was 43000
$ws.Range("J131").Value = 45000  # was 43000
$ws.Range("L131").Value = 45000  # was 43000
$ws.Range("N131").Value = -55080  # was -53080
$ws.Range("H139").Value = 59350  # was 59800
$ws.Range("J139").Value = 59350  # was 59800
$ws.Range("L139").Value = 59350  # was 59800
$ws.Range("N139").Value = -70080  # was -70840
